# Restore the "Rules" sheet's C10 cell value to 1 (was 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
